$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 7220
$ws.Range("E2").Value = 453
$ws.Range("F2").Value = 453
$ws.Range("G2").Value = 561
$ws.Range("H2").Value = 448
$ws.Range("I2").Value = 454
$ws.Range("J2").Value = -6
$ws.Range("K2").Value = 9716
$ws.Range("L2").Value = 3213
$ws.Range("M2").Value = 6503
$ws.Range("N2").Value = 6415
$ws.Range("O2").Value = 88
$ws.Range("P2").Value = 119
$ws.Range("Q2").Value = 595
$ws.Range("R2").Value = -768
$ws.Range("S2").Value = 339
$ws.Range("T2").Value = 606
$ws.Range("U2").Value = -11
$ws.Range("V2").Value = 1649
$ws.Range("W2").Value = 6.27
$ws.Range("X2").Value = 6.21
$ws.Range("Y2").Value = 7.34
$ws.Range("Z2").Value = 5.07
$ws.Range("AA2").Value = 49.41
$ws.Range("AB2").Value = 5356.94
$ws.Range("AC2").Value = 1915
$ws.Range("AD2").Value = 14.38
$ws.Range("AE2").Value = 30954
$ws.Range("AF2").Value = 0.89
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 0.54
$ws.Range("AI2").Value = 6.84
$ws.Range("AJ2").Value = 23727020

# Row 3
$ws.Range("D3").Value = 9105
$ws.Range("E3").Value = 849
$ws.Range("F3").Value = 849
$ws.Range("G3").Value = 935
$ws.Range("H3").Value = 747
$ws.Range("I3").Value = 727
$ws.Range("J3").Value = 19
$ws.Range("K3").Value = 10645
$ws.Range("L3").Value = 3391
$ws.Range("M3").Value = 7254
$ws.Range("N3").Value = 7150
$ws.Range("O3").Value = 104
$ws.Range("P3").Value = 119
$ws.Range("Q3").Value = 1147
$ws.Range("R3").Value = -1227
$ws.Range("S3").Value = -52
$ws.Range("T3").Value = 418
$ws.Range("U3").Value = 729
$ws.Range("V3").Value = 1614
$ws.Range("W3").Value = 9.33
$ws.Range("X3").Value = 8.2
$ws.Range("Y3").Value = 10.72
$ws.Range("Z3").Value = 7.33
$ws.Range("AA3").Value = 46.74
$ws.Range("AB3").Value = 5946.31
$ws.Range("AC3").Value = 3066
$ws.Range("AD3").Value = 11.17
$ws.Range("AE3").Value = 34504
$ws.Range("AF3").Value = 0.99
$ws.Range("AG3").Value = 250
$ws.Range("AH3").Value = 0.73
$ws.Range("AI3").Value = 7.12
$ws.Range("AJ3").Value = 23727020

# Row 4
$ws.Range("D4").Value = 7957
$ws.Range("E4").Value = 652
$ws.Range("F4").Value = 922
$ws.Range("G4").Value = 785
$ws.Range("H4").Value = 840
$ws.Range("I4").Value = 800
$ws.Range("J4").Value = 40
$ws.Range("K4").Value = 11155
$ws.Range("L4").Value = 3251
$ws.Range("M4").Value = 7905
$ws.Range("N4").Value = 7767
$ws.Range("O4").Value = 138
$ws.Range("P4").Value = 119
$ws.Range("Q4").Value = 1361
$ws.Range("R4").Value = -739
$ws.Range("S4").Value = -230
$ws.Range("T4").Value = 486
$ws.Range("U4").Value = 875
$ws.Range("V4").Value = 1432
$ws.Range("W4").Value = 8.2
$ws.Range("X4").Value = 10.55
$ws.Range("Y4").Value = 10.73
$ws.Range("Z4").Value = 7.7
$ws.Range("AA4").Value = 41.12
$ws.Range("AB4").Value = 6564.42
$ws.Range("AC4").Value = 3372
$ws.Range("AD4").Value = 10.62
$ws.Range("AE4").Value = 37478
$ws.Range("AF4").Value = 0.96
$ws.Range("AG4").Value = 350
$ws.Range("AH4").Value = 0.98
$ws.Range("AI4").Value = 9.06
$ws.Range("AJ4").Value = 23727020

# Row 5
$ws.Range("D5").Value = 8779
$ws.Range("E5").Value = 595
$ws.Range("F5").Value = 595
$ws.Range("G5").Value = 510
$ws.Range("H5").Value = 535
$ws.Range("I5").Value = 494
$ws.Range("J5").Value = 41
$ws.Range("K5").Value = 7786
$ws.Range("L5").Value = 2658
$ws.Range("M5").Value = 5128
$ws.Range("N5").Value = 4941
$ws.Range("O5").Value = 187
$ws.Range("P5").Value = 50
$ws.Range("Q5").Value = 798
$ws.Range("R5").Value = -451
$ws.Range("S5").Value = -509
$ws.Range("T5").Value = 402
$ws.Range("U5").Value = 395
$ws.Range("V5").Value = 923
$ws.Range("W5").Value = 6.78
$ws.Range("X5").Value = 6.09
$ws.Range("Y5").Value = 7.77
$ws.Range("Z5").Value = 5.65
$ws.Range("AA5").Value = 51.84
$ws.Range("AB5").Value = 15751.13
$ws.Range("AC5").Value = 2199
$ws.Range("AD5").Value = 14.53
$ws.Range("AE5").Value = 56602
$ws.Range("AF5").Value = 0.56
$ws.Range("AG5").Value = 700
$ws.Range("AH5").Value = 2.19
$ws.Range("AI5").Value = 1.24
$ws.Range("AJ5").Value = 9994005

# Row 6
$ws.Range("D6").Value = 9085
$ws.Range("E6").Value = 717
$ws.Range("F6").Value = 717
$ws.Range("G6").Value = 788
$ws.Range("H6").Value = 585
$ws.Range("I6").Value = 558
$ws.Range("K6").Value = 8708
$ws.Range("L6").Value = 3072
$ws.Range("M6").Value = 5635
$ws.Range("N6").Value = 5425
$ws.Range("P6").Value = 50
$ws.Range("Q6").Value = 742
$ws.Range("R6").Value = -494
$ws.Range("S6").Value = 363
$ws.Range("T6").Value = 466
$ws.Range("U6").Value = 276
$ws.Range("V6").Value = 1346
$ws.Range("W6").Value = 7.89
$ws.Range("X6").Value = 6.44
$ws.Range("Y6").Value = 10.76
$ws.Range("Z6").Value = 7.09
$ws.Range("AA6").Value = 54.52
$ws.Range("AB6").Value = 16701.48
$ws.Range("AC6").Value = 5581
$ws.Range("AD6").Value = 5.73
$ws.Range("AE6").Value = 62158
$ws.Range("AF6").Value = 0.51
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 3.13
$ws.Range("AI6").Value = 15.65
$ws.Range("AJ6").Value = 9994005

# Rows 7-9: clear all data columns, keep only A, B, C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
